# Generate Report for Handback
#
# For the "65d94e59-5e31-4c4a-97fb-5e441f88bd68" source file, a handback
# xliff was found/processed for both the zh-cn and de-de targets. This
# records the discovered handback target file + datetime, and flags an
# error because the handback file version isn't the latest one available
# (current vs. latest commit of the underlying .md file differ).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fa56f608d1cd574ad662fa34717f0292661a74f/e2e/65d94e59-5e31-4c4a-97fb-5e441f88bd68.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6db441a12dea4ee74b30cce6c25b57c4a8ecd87/e2e/65d94e59-5e31-4c4a-97fb-5e441f88bd68.md."
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fa56f608d1cd574ad662fa34717f0292661a74f/e2e/65d94e59-5e31-4c4a-97fb-5e441f88bd68.md"
$displayName = "65d94e59-5e31-4c4a-97fb-5e441f88bd68.md"

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$TargetFile,
        [string]$HandbackDate
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # I5: Latest Target File - link to the md file for the current commit
    $ws.Range("I5").Value = $displayName
    $ws.Hyperlinks.Add($ws.Range("I5"), $currentUrl, "", "", $displayName) | Out-Null

    # J5: Latest Handback File
    $ws.Range("J5").Value = $TargetFile

    # K5: Latest Handback DateTime
    $ws.Range("K5").Value = $HandbackDate

    # P5: Error Detail
    $ws.Range("P5").Value = $errorDetail

    # Widen the columns that now hold the newly populated / long text so the
    # report is readable.
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(16).ColumnWidth = 40
}

Update-HandbackRow "zh-cn" "65d94e59-5e31-4c4a-97fb-5e441f88bd68.27a3c81b24a7edb78fd15c28029d67d6875747f8.zh-cn.xlf" "2016-09-09 11:58:13"
Update-HandbackRow "de-de" "65d94e59-5e31-4c4a-97fb-5e441f88bd68.27a3c81b24a7edb78fd15c28029d67d6875747f8.de-de.xlf" "2016-09-09 11:58:32"
